$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1606425702811245
$ws.Range("C2").Value = 0.6104417670682731
$ws.Range("J2").Value = 0.02409638554216868
$ws.Range("P2").Value = 0.1124497991967871
$ws.Range("S2").Value = 0.09236947791164658
$ws.Range("B3").Value = 0.01886792452830189
$ws.Range("C3").Value = 0.02515723270440252
$ws.Range("J3").Value = 0.02515723270440252
$ws.Range("P3").Value = 0.7924528301886793
$ws.Range("S3").Value = 0.1383647798742138
$ws.Range("J4").Value = 0.07317073170731707
$ws.Range("P4").Value = 0.7317073170731707
$ws.Range("S4").Value = 0.1951219512195122
$ws.Range("B6").Value = 0.03773584905660377
$ws.Range("D6").Value = 0.009433962264150943
$ws.Range("F6").Value = 0.05188679245283019
$ws.Range("J6").Value = 0.2735849056603774
$ws.Range("O6").Value = 0.01886792452830189
$ws.Range("Q6").Value = 0.1462264150943396
$ws.Range("R6").Value = 0.05660377358490566
$ws.Range("S6").Value = 0.4056603773584906
$ws.Range("B7").Value = 0.09844559585492228
$ws.Range("D7").Value = 0.02072538860103627
$ws.Range("E7").Value = 0.005181347150259068
$ws.Range("F7").Value = 0.05699481865284974
$ws.Range("J7").Value = 0.1606217616580311
$ws.Range("O7").Value = 0.0310880829015544
$ws.Range("Q7").Value = 0.1295336787564767
$ws.Range("R7").Value = 0.1036269430051813
$ws.Range("S7").Value = 0.3937823834196891
$ws.Range("B8").Value = 0.08053691275167785
$ws.Range("D8").Value = 0.02013422818791946
$ws.Range("F8").Value = 0.07606263982102908
$ws.Range("J8").Value = 0.116331096196868
$ws.Range("O8").Value = 0.01789709172259508
$ws.Range("Q8").Value = 0.1901565995525727
$ws.Range("R8").Value = 0.08053691275167785
$ws.Range("S8").Value = 0.4183445190156599
$ws.Range("B9").Value = 0.07926829268292683
$ws.Range("D9").Value = 0.01829268292682927
$ws.Range("F9").Value = 0.04878048780487805
$ws.Range("J9").Value = 0.1158536585365854
$ws.Range("O9").Value = 0.03658536585365853
$ws.Range("Q9").Value = 0.2073170731707317
$ws.Range("R9").Value = 0.0975609756097561
$ws.Range("S9").Value = 0.3963414634146342
$ws.Range("B10").Value = 0.1142857142857143
$ws.Range("D10").Value = 0.02232142857142857
$ws.Range("F10").Value = 0.08839285714285715
$ws.Range("J10").Value = 0.1160714285714286
$ws.Range("O10").Value = 0.01964285714285714
$ws.Range("Q10").Value = 0.2017857142857143
$ws.Range("R10").Value = 0.07232142857142858
$ws.Range("S10").Value = 0.3651785714285714
$ws.Range("G11").Value = 0.1348314606741573
$ws.Range("J11").Value = 0.09737827715355805
$ws.Range("K11").Value = 0.2059925093632959
$ws.Range("L11").Value = 0.550561797752809
$ws.Range("S11").Value = 0.01123595505617977
$ws.Range("G12").Value = 0.7880794701986755
$ws.Range("J12").Value = 0.1854304635761589
$ws.Range("L12").Value = 0.01324503311258278
$ws.Range("S12").Value = 0.01324503311258278
$ws.Range("G13").Value = 0.7543859649122807
$ws.Range("J13").Value = 0.2280701754385965
$ws.Range("S13").Value = 0.01754385964912281
$ws.Range("F14").Value = 0.5
$ws.Range("J14").Value = 0.5
$ws.Range("F15").Value = 0.004784688995215311
$ws.Range("H15").Value = 0.1674641148325359
$ws.Range("I15").Value = 0.04784688995215311
$ws.Range("J15").Value = 0.3971291866028708
$ws.Range("K15").Value = 0.0430622009569378
$ws.Range("M15").Value = 0.004784688995215311
$ws.Range("N15").Value = 0.004784688995215311
$ws.Range("O15").Value = 0.03827751196172249
$ws.Range("S15").Value = 0.291866028708134
$ws.Range("F16").Value = 0.03389830508474576
$ws.Range("H16").Value = 0.2033898305084746
$ws.Range("I16").Value = 0.07344632768361582
$ws.Range("J16").Value = 0.3615819209039548
$ws.Range("K16").Value = 0.1016949152542373
$ws.Range("M16").Value = 0.05084745762711865
$ws.Range("O16").Value = 0.05649717514124294
$ws.Range("S16").Value = 0.1186440677966102
$ws.Range("F17").Value = 0.01256281407035176
$ws.Range("H17").Value = 0.1683417085427136
$ws.Range("I17").Value = 0.1030150753768844
$ws.Range("J17").Value = 0.4296482412060301
$ws.Range("K17").Value = 0.08542713567839195
$ws.Range("M17").Value = 0.01758793969849246
$ws.Range("O17").Value = 0.05778894472361809
$ws.Range("S17").Value = 0.1256281407035176
$ws.Range("F18").Value = 0.01197604790419162
$ws.Range("H18").Value = 0.2335329341317365
$ws.Range("I18").Value = 0.07784431137724551
$ws.Range("J18").Value = 0.407185628742515
$ws.Range("K18").Value = 0.1017964071856287
$ws.Range("M18").Value = 0.005988023952095809
$ws.Range("O18").Value = 0.0658682634730539
$ws.Range("S18").Value = 0.09580838323353294
$ws.Range("F19").Value = 0.01543739279588336
$ws.Range("H19").Value = 0.2341337907375643
$ws.Range("I19").Value = 0.07632933104631218
$ws.Range("J19").Value = 0.3327615780445969
$ws.Range("K19").Value = 0.114065180102916
$ws.Range("M19").Value = 0.03173241852487135
$ws.Range("N19").Value = 0.0008576329331046312
$ws.Range("O19").Value = 0.06775300171526587
$ws.Range("S19").Value = 0.1269296740994854
